$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repulled / recalculated data
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = -5
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = -4
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = -3
$ws.Range("F9").Value = -3
$ws.Range("F10").Value = 1
$ws.Range("F11").Value = -6
$ws.Range("F12").Value = 9
$ws.Range("F15").Value = -2
$ws.Range("F18").Value = -3
$ws.Range("F19").Value = -1
$ws.Range("F20").Value = 3
$ws.Range("F21").Value = -2
$ws.Range("F22").Value = -1
